$d = $word.ActiveDocument

# --- Step 1: drop the _GoBack bookmark from its current spot (right before
# "the data that we got back from Open Weather Map."). The surrounding runs
# are left completely untouched, so removing the (empty/zero-width) bookmark
# simply merges that sentence back into the flow of the paragraph with no
# bookmark splitting it anymore.
$bookmark = $d.Bookmarks.Item("_GoBack")
$bookmark.Delete()

# --- Step 2: find + delete the paragraph that holds the legacy ActiveX/OLE
# control (the one right after "So all that and more, I'll see you on the
# next lesson."). Deleting the whole paragraph range (including its mark)
# removes the object/control run completely and folds the following
# (empty) paragraph up into its place.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*So all that and more, I'll see you on the next lesson*") {
        $targetIndex = $i + 1
        break
    }
}
if ($targetIndex -ne -1) {
    $objPara = $d.Paragraphs.Item($targetIndex)
    $objPara.Range.Delete()
}

# --- Step 3: put a fresh, empty "_GoBack" bookmark into the (now empty)
# final paragraph of the document. Word refuses to stash a zero-width
# bookmark exactly at a paragraph mark in one shot, so we briefly insert a
# placeholder character, wrap the bookmark around it, then delete the
# placeholder again -- the bookmark collapses back to zero width and stays
# put, sitting alone in its paragraph.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastStart = $lastPara.Range.Start

$placeholder = $d.Range($lastStart, $lastStart)
$placeholder.InsertBefore("X")

$wrap = $d.Range($lastStart, $lastStart + 1)
$d.Bookmarks.Add("_GoBack", $wrap)

$cleanup = $d.Range($lastStart, $lastStart + 1)
$cleanup.Delete()
